# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (want-to-go count) figures in the F column of the
# "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 7801
$ws1.Range("F6").Value  = 570
$ws1.Range("F7").Value  = 1186
$ws1.Range("F8").Value  = 211
$ws1.Range("F10").Value = 173

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 7801
$ws4.Range("F6").Value  = 570
$ws4.Range("F7").Value  = 1186
$ws4.Range("F8").Value  = 211
$ws4.Range("F11").Value = 173
